$d = $word.ActiveDocument

$replacements = @(
    @("825÷3=", "733÷7="),
    @("258÷8=", "742÷6="),
    @("205÷8=", "342÷5="),
    @("981÷6=", "622÷5="),
    @("318÷3=", "557÷7="),
    @("836÷4=", "400÷4="),
    @("432÷6=", "527÷7="),
    @("908÷3=", "456÷8="),
    @("625÷4=", "328÷7="),
    @("676÷3=", "233÷6="),
    @("190÷9=", "188÷3="),
    @("881÷7=", "925÷2="),
    @("937÷5=", "296÷9="),
    @("294÷2=", "465÷7="),
    @("143÷7=", "747÷8="),
    @("842÷9=", "705÷8="),
    @("741÷4=", "491÷9="),
    @("590÷4=", "920÷5="),
    @("106÷4=", "122÷9="),
    @("986÷9=", "323÷3="),
    @("757÷6=", "707÷9="),
    @("522÷9=", "268÷2="),
    @("137÷7=", "733÷2="),
    @("518÷8=", "147÷9="),
    @("101÷3=", "342÷3=")
)

foreach ($pair in $replacements) {
    $old = $pair[0]
    $new = $pair[1]
    $d.Content.Find.Execute($old, $true, $false, $false, $false, $false,
                             $true, 1, $false, $new, 2)
}

Write-Output "Done: $($replacements.Count) replacements applied"
